$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 5: update the change-log cell (C5) to append the new bullet point ---
$ws.Range("C5").Value = "-Agrupation and desagrupation to be done in functions.`n-Change reproduction and distribution to two parts.`n-UI: Delete rows according to working functionality.`n-UI: condicionate IF to be associated.`n-Implement mutations.`n-Implement save and load configurations.`n-Document every function.`n-The program does not work with 2 or less niches.`n-When an actor dies, the recipient should deassociate.`n-Order for SG shall have into account if the specie is associated or not.`n-UI: automatically fill aggrupation data.`n-SI is not doing the random selection properly."

# --- Row 6: new version entry 0.1.3 ---
$ws.Range("A6").Value = "0.1.3"
$ws.Range("B6").Value = "AUTOMATA CELULAR - copia (10)"
$ws.Range("C6").Value = $ws.Range("C5").Value
$ws.Range("D6").Value = "-SI random selection corrected."
$ws.Range("E6").Value = "Python 3.6.1"
$ws.Range("F6").Value = "Qt version: 5.6.2`nSIP version: 4.18`nPyQt version: 5.6"
$ws.Range("G6").Value = " PyInstaller 3.3.1"

# --- Row heights ---
$ws.Range("A5:G5").RowHeight = 187.2
$ws.Range("A6:G6").RowHeight = 187.2

# --- Column width ---
$ws.Range("B1").ColumnWidth = 28.5546875

# --- View / selection ---
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("H6").Select()
